$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "2.0-harinita,1.0-huevos,1.0-vainilla,5.0-leche,"
$ws.Range("C3").Value = "1.0-harinita,5.0-huevos,2.0-manzana,"
$ws.Range("C4").Value = "5.0-harinita,2.0-huevos,1.0-vainilla,"
$ws.Range("C5").Value = "5.0-harinita,5.0-huevos,"
$ws.Range("C6").Value = "4.0-harinita,1.0-crema,5.0-merengue,2.0-limon,5.0-huevos,"
